# Update column B (CC) values on the active worksheet ("Main Dashboard" data)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5444.512500000001
$ws.Range("B3").Value = 5387.3435
$ws.Range("B4").Value = 5357.2225
$ws.Range("B5").Value = 5289.154500000001
$ws.Range("B6").Value = 5361.8285
$ws.Range("B7").Value = 5442.6015
$ws.Range("B8").Value = 5816.530999999999
$ws.Range("B9").Value = 7034.6115
$ws.Range("B10").Value = 8596.6895
$ws.Range("B11").Value = 13942.873
$ws.Range("B12").Value = 15769.3095
$ws.Range("B13").Value = 15209.229
$ws.Range("B14").Value = 15159.795
$ws.Range("B15").Value = 15449.4445
$ws.Range("B16").Value = 15588.433
$ws.Range("B17").Value = 16020.6585
$ws.Range("B18").Value = 16262.3825
$ws.Range("B19").Value = 15862.609
$ws.Range("B20").Value = 14887.418
$ws.Range("B21").Value = 13341.0375
$ws.Range("B22").Value = 11837.672
$ws.Range("B23").Value = 9490.991999999998
$ws.Range("B24").Value = 6891.794
$ws.Range("B25").Value = 5770.9715
